$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 520.46875
$ws.Range("J33").Value = 1687.8
$ws.Range("L33").Value = 1687.8
$ws.Range("N33").Value = -2145.8
$ws.Range("H116").Value = 1500
$ws.Range("I116").Value = 1500
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1500
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1942
$ws.Range("N116").ClearContents()
$ws.Range("H121").Value = 661.5714
$ws.Range("J121").Value = 821
$ws.Range("L121").Value = 2463
$ws.Range("N121").Value = -5957

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4133921.5
$ws.Range("I45").Value = 4547138.5
$ws.Range("J45").Value = 1750
$ws.Range("K45").Value = 4547138.5
$ws.Range("L45").Value = 1750
$ws.Range("M45").Value = -4546761.5
$ws.Range("N45").Value = -2504
$ws.Range("H110").Value = 856.6842
$ws.Range("I110").Value = 832.05554
$ws.Range("K110").Value = 832.05554
$ws.Range("M110").Value = 1212.94446

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 87019.86
$ws.Range("J81").Value = 87019.86
$ws.Range("L81").Value = 87019.86
$ws.Range("N81").Value = -89141.86
$ws.Range("H84").Value = 87019.86
$ws.Range("J84").Value = 87019.86
$ws.Range("L84").Value = 261059.58
$ws.Range("N84").Value = -271667.58
$ws.Range("H134").Value = 838208.0600000001
$ws.Range("I134").Value = 1337762
$ws.Range("J134").Value = 5618.222
$ws.Range("K134").Value = 4013286
$ws.Range("L134").Value = 16854.666
$ws.Range("M134").Value = -4010751
$ws.Range("N134").Value = -21924.666
$ws.Range("H135").Value = 35857.145
$ws.Range("J135").Value = 35857.145
$ws.Range("L135").Value = 35857.145
$ws.Range("N135").Value = -45997.145

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2606.7778
$ws.Range("I31").Value = 2112.0908
$ws.Range("J31").Value = 3384.1428
$ws.Range("K31").Value = 2112.0908
$ws.Range("L31").Value = 3384.1428
$ws.Range("M31").Value = -1817.0908
$ws.Range("N31").Value = -3974.1428
$ws.Range("H34").Value = 2606.7778
$ws.Range("I34").Value = 2112.0908
$ws.Range("J34").Value = 3384.1428
$ws.Range("K34").Value = 2112.0908
$ws.Range("L34").Value = 3384.1428
$ws.Range("M34").Value = -1910.0908
$ws.Range("N34").Value = -3788.1428
$ws.Range("H141").Value = 41971
$ws.Range("J141").Value = 41971
$ws.Range("L141").Value = 41971
$ws.Range("N141").Value = -52331

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 809.36365
$ws.Range("I68").Value = 700
$ws.Range("J68").Value = 820.3
$ws.Range("K68").Value = 2100
$ws.Range("L68").Value = 2460.9
$ws.Range("M68").Value = -1289
$ws.Range("N68").Value = -4082.9
$ws.Range("H71").Value = 809.36365
$ws.Range("I71").Value = 700
$ws.Range("J71").Value = 820.3
$ws.Range("K71").Value = 6300
$ws.Range("L71").Value = 7382.7
$ws.Range("M71").Value = -2244
$ws.Range("N71").Value = -15494.7
$ws.Range("H131").Value = 926.48
$ws.Range("J131").Value = 957.95746
$ws.Range("L131").Value = 2873.87238
$ws.Range("N131").Value = -12953.87238
$ws.Range("H134").Value = 2956.3333
$ws.Range("I134").Value = 2395.926
$ws.Range("J134").Value = 8000
$ws.Range("K134").Value = 7187.778
$ws.Range("L134").Value = 24000
$ws.Range("M134").Value = -2117.778
$ws.Range("N134").Value = -34140

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2338.9092
$ws.Range("I80").Value = 2331.9092
$ws.Range("J80").Value = 2345.9092
$ws.Range("K80").Value = 2331.9092
$ws.Range("L80").Value = 2345.9092
$ws.Range("M80").Value = -1333.9092
$ws.Range("N80").Value = -4341.9092
$ws.Range("H83").Value = 2338.9092
$ws.Range("I83").Value = 2331.9092
$ws.Range("J83").Value = 2345.9092
$ws.Range("K83").Value = 11659.546
$ws.Range("L83").Value = 11729.546
$ws.Range("M83").Value = -6667.546
$ws.Range("N83").Value = -21713.546
$ws.Range("H102").Value = 13890216
$ws.Range("I102").Value = 27778456
$ws.Range("J102").Value = 1975
$ws.Range("K102").Value = 27778456
$ws.Range("L102").Value = 1975
$ws.Range("M102").Value = -27776834
$ws.Range("N102").Value = -5219
$ws.Range("H126").Value = 41668444
$ws.Range("I126").Value = 83334264
$ws.Range("J126").Value = 2628.5
$ws.Range("K126").Value = 250002792
$ws.Range("L126").Value = 7885.5
$ws.Range("M126").Value = -250000322
$ws.Range("N126").Value = -12825.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H81").Value = 30000
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 30000
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 30000
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -31996
$ws.Range("H82").Value = 987.5
$ws.Range("I82").Value = 975
$ws.Range("J82").Value = 990
$ws.Range("K82").Value = 975
$ws.Range("L82").Value = 990
$ws.Range("M82").Value = -614
$ws.Range("N82").Value = -1712
$ws.Range("H84").Value = 30000
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 30000
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 90000
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -99984
$ws.Range("H85").Value = 987.5
$ws.Range("I85").Value = 975
$ws.Range("J85").Value = 990
$ws.Range("K85").Value = 975
$ws.Range("L85").Value = 990
$ws.Range("M85").Value = 273
$ws.Range("N85").Value = -3486
$ws.Range("H122").Value = 3053.2
$ws.Range("I122").Value = 1899
$ws.Range("J122").Value = 3230.7693
$ws.Range("K122").Value = 5697
$ws.Range("L122").Value = 9692.3079
$ws.Range("M122").Value = -3247
$ws.Range("N122").Value = -14592.3079

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 20383.188
$ws.Range("J109").Value = 20383.188
$ws.Range("L109").Value = 20383.188
$ws.Range("N109").Value = -23157.188
$ws.Range("H113").Value = 444.70834
$ws.Range("I113").Value = 330
$ws.Range("J113").Value = 880.6
$ws.Range("K113").Value = 990
$ws.Range("L113").Value = 2641.8
$ws.Range("M113").Value = 1180
$ws.Range("N113").Value = -6981.8
$ws.Range("H115").Value = 30000
$ws.Range("J115").Value = 30000
$ws.Range("L115").Value = 30000
$ws.Range("N115").Value = -33134
$ws.Range("H118").Value = 29329.666
$ws.Range("J118").Value = 29000
$ws.Range("L118").Value = 29000
$ws.Range("N118").Value = -32314
